$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 250
$ws.Range("F4").Value = 844
$ws.Range("F6").Value = 416
$ws.Range("F7").Value = 604
$ws.Range("F8").Value = 223
$ws.Range("F10").Value = 353
$ws.Range("F11").Value = 153
$ws.Range("F12").Value = 689
$ws.Range("F13").Value = 91
$ws.Range("F14").Value = 1823
$ws.Range("F15").Value = 363
$ws.Range("F16").Value = 3592
$ws.Range("F17").Value = 331
$ws.Range("F18").Value = 496
$ws.Range("F19").Value = 5
$ws.Range("F20").Value = 58
$ws.Range("F21").Value = 145

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 223
$ws.Range("F5").Value = 21
$ws.Range("F6").Value = 119
$ws.Range("F7").Value = 479
$ws.Range("F13").Value = 94
$ws.Range("F21").Value = 22

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5357
$ws.Range("F3").Value = 329
$ws.Range("F4").Value = 279

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5357
$ws.Range("F4").Value = 329
$ws.Range("F6").Value = 279
$ws.Range("F7").Value = 250
$ws.Range("F8").Value = 223
$ws.Range("F10").Value = 21
$ws.Range("F11").Value = 119
$ws.Range("F12").Value = 479
$ws.Range("F13").Value = 479
$ws.Range("F14").Value = 844
$ws.Range("F18").Value = 416
$ws.Range("F19").Value = 604
$ws.Range("F20").Value = 223
$ws.Range("F23").Value = 353
$ws.Range("F24").Value = 153
$ws.Range("F27").Value = 689
$ws.Range("F28").Value = 91
$ws.Range("F29").Value = 94
$ws.Range("F30").Value = 1823
$ws.Range("F31").Value = 363
$ws.Range("F32").Value = 3593
$ws.Range("F34").Value = 331
$ws.Range("F35").Value = 496
$ws.Range("F36").Value = 5
$ws.Range("F37").Value = 58
$ws.Range("F39").Value = 145
$ws.Range("F46").Value = 22
